# Update the "Date" and "FHIR Version" metadata values, and a few content
# corrections on the Elements sheet, to match the new IG build.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$wsMeta.Range("B15").Value = "4.0.1"

$wsElements = $wb.Worksheets.Item("Elements")

# Extension root constraint text (drop the "unless an empty Parameters resource" clause)
$wsElements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.id Type(s) corrected from "id" to "string"
$wsElements.Range("K3").Value = "string`n"

# Extension.value[x] short description: R4B -> R4 doc link
$wsElements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
